# Generate Report for Handoff
# Refreshes the localization-status report with a new handoff run:
#  - the previous "transform failed" row (95b843ef...) is gone
#  - a new source file (410c6b89-7040-4a4a-8bea-1c625417d42d.md) was
#    handed off, producing new xlf packages with fresh timestamps
#  - the ".localization-config" row moves up to take the old row's place

$wb = $excel.ActiveWorkbook

$oldUuid = "5033ca86-d84c-4b42-80e6-47992eb0236e"
$newUuid = "410c6b89-7040-4a4a-8bea-1c625417d42d"
$oldHash = "74d46d7087bef962109ba6b874c1e2affd3a8cd3"
$newHash = "dc9f4ca938b89784312d30588b80f361d8a0c0d0"

$newMdName = "$newUuid.md"
$cfgUrl = "https://github.com/OpenLocalizationTest/oltest/blob/c41e2daa0e0f98f6f121329e0013e61b33d66eba/.localization-config"

# ------------------------------------------------------------------
# Sheet "Overview"
# ------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Overview")

# Drop the "Handoff transform failed" row; this shifts the
# ".localization-config" row up one and keeps its existing styling.
$ws1.Rows.Item(3).Delete()

# New source file name for the handed-off document.
$ws1.Range("A2").Value = $newMdName

$mdUrl1 = "https://github.com/OpenLocalizationTest/oltest/blob/c41e2daa0e0f98f6f121329e0013e61b33d66eba/e2e/$newMdName"

$ws1.Hyperlinks.Delete()
$ws1.Hyperlinks.Add($ws1.Range("A2"), $mdUrl1, [Type]::Missing, [Type]::Missing, $newMdName) | Out-Null
$ws1.Hyperlinks.Add($ws1.Range("A3"), $cfgUrl, [Type]::Missing, [Type]::Missing, ".localization-config") | Out-Null

# ------------------------------------------------------------------
# Sheet "zh-cn"
# ------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("zh-cn")

$ws2.Rows.Item(3).Delete()

$ws2.Range("A2").Value = $newMdName
$newXlf2 = "$newUuid.$newHash.zh-cn.xlf"
$ws2.Range("C2").Value = $newXlf2
$ws2.Range("D2").Value = "2016-02-18 10:11:57"

$mdUrl2 = "https://github.com/OpenLocalizationTest/oltest/blob/c41e2daa0e0f98f6f121329e0013e61b33d66eba/e2e/$newMdName"
$xlfUrl2 = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/8ebc366fe557f818276319f8305b7934e840e9b2/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/$newXlf2"

$ws2.Hyperlinks.Delete()
$ws2.Hyperlinks.Add($ws2.Range("A2"), $mdUrl2, [Type]::Missing, [Type]::Missing, $newMdName) | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("C2"), $xlfUrl2, [Type]::Missing, [Type]::Missing, $newXlf2) | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("A3"), $cfgUrl, [Type]::Missing, [Type]::Missing, ".localization-config") | Out-Null

# ------------------------------------------------------------------
# Sheet "de-de"
# ------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("de-de")

$ws3.Rows.Item(3).Delete()

$ws3.Range("A2").Value = $newMdName
$newXlf3 = "$newUuid.$newHash.de-de.xlf"
$ws3.Range("C2").Value = $newXlf3
$ws3.Range("D2").Value = "2016-02-18 10:12:08"

$mdUrl3 = "https://github.com/OpenLocalizationTest/oltest/blob/c41e2daa0e0f98f6f121329e0013e61b33d66eba/e2e/$newMdName"
$xlfUrl3 = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/cc17901c8175cee06bfe92ed02c8416c44def831/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/$newXlf3"

$ws3.Hyperlinks.Delete()
$ws3.Hyperlinks.Add($ws3.Range("A2"), $mdUrl3, [Type]::Missing, [Type]::Missing, $newMdName) | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("C2"), $xlfUrl3, [Type]::Missing, [Type]::Missing, $newXlf3) | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("A3"), $cfgUrl, [Type]::Missing, [Type]::Missing, ".localization-config") | Out-Null

Write-Host "Report regenerated for handoff."
